$d = $word.ActiveDocument

$find = ".\RandomMapper.py | python .\RandomReducer.py"
$replace = ".\RandomMapper.py | Sort-Object | python .\RandomReducer.py"

$range = $d.Content
$range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
